# Add new test-case data (verifyLoginWithValidCred) to the "testData" sheet,
# including a hyperlink cell pointing at the CRM login page.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testData")

$ws.Range("A11").Value = "verifyLoginWithValidCred"

$ws.Range("A12").Value = "username"
$ws.Range("B12").Value = "password"

$ws.Range("A13").Value = "w2ajava@way2automation.com"
$ws.Range("B13").Value = "Tcs@12345"

$null = $ws.Hyperlinks.Add($ws.Range("B13"), "http://www.way2automation.com/way2auto_jquery/crm/")

$null = $ws.Range("B13").Select()
